{"js": "// The answers in the practice table (\"two-digit number divided by\n// one-digit number\") were regenerated; every populated cell's text is\n// replaced in place (row/column position), preserving all paragraph\n// and run formatting (font, size, alignment) already on the cell.\n//\n// Because a couple of old/new values repeat (e.g. \"79\u00f73=26, 1\" occurs\n// twice before the edit, and \"68\u00f76=11, 2\" / \"77\u00f73=25, 2\" each occur\n// twice after), matching must be done by (row, column) position, NOT\n// by searching for the old text.\n\nconst newValues = [\n  [\"68\u00f76=11, 2\", \"18\u00f77=2, 4\", \"86\u00f73=28, 2\", \"86\u00f79=9, 5\", \"51\u00f73=17, 0\"],\n  [\"32\u00f78=4, 0\", \"58\u00f76=9, 4\", \"55\u00f72=27, 1\", \"21\u00f75=4, 1\", \"63\u00f76=10, 3\"],\n  [\"25\u00f72=12, 1\", \"77\u00f73=25, 2\", \"77\u00f73=25, 2\", \"53\u00f74=13, 1\", \"80\u00f77=11, 3\"],\n  [\"48\u00f73=16, 0\", \"55\u00f75=11, 0\", \"14\u00f79=1, 5\", \"17\u00f74=4, 1\", \"60\u00f73=20, 0\"],\n  [\"30\u00f75=6, 0\", \"68\u00f76=11, 2\", \"61\u00f72=30, 1\", \"79\u00f74=19, 3\", \"69\u00f74=17, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Only every 4th row (0, 4, 8, 12, 16) actually holds the answer text;\n// the rows between them are spacer rows with empty paragraphs. Walk\n// the rows that have cell content and apply the next answer set, so\n// the mapping is robust even if spacer-row geometry differs.\nlet dataRowIndex = 0;\nfor (const row of rows) {\n  const cells = row.cells.items;\n  if (dataRowIndex >= newValues.length) break;\n  const rowValues = newValues[dataRowIndex];\n\n  // Detect whether this row is a data row (has any non-empty text).\n  for (const cell of cells) {\n    cell.load(\"value\");\n  }\n  await context.sync();\n\n  const hasText = cells.some((c) => (c.value || \"\").trim().length > 0);\n  if (!hasText) continue;\n\n  for (let c = 0; c < cells.length && c < rowValues.length; c++) {\n    cells[c].value = rowValues[c];\n  }\n  dataRowIndex++;\n}\n\nawait context.sync();\n", "ps1": "# The answers in the practice table (\"two-digit number divided by\n# one-digit number\") were regenerated; every populated cell's text is\n# replaced in place (row/column position), preserving all paragraph\n# and run formatting (font, size, alignment) already on the cell.\n#\n# Because a couple of old/new values repeat (e.g. \"79\u00f73=26, 1\" occurs\n# twice before the edit, and \"68\u00f76=11, 2\" / \"77\u00f73=25, 2\" each occur\n# twice after), matching is done by (row, column) position via\n# Table.Cell(r, c), NOT by searching/replacing the old text.\n\n$newValues = @(\n    ,@(\"68\u00f76=11, 2\", \"18\u00f77=2, 4\", \"86\u00f73=28, 2\", \"86\u00f79=9, 5\", \"51\u00f73=17, 0\")\n    ,@(\"32\u00f78=4, 0\", \"58\u00f76=9, 4\", \"55\u00f72=27, 1\", \"21\u00f75=4, 1\", \"63\u00f76=10, 3\")\n    ,@(\"25\u00f72=12, 1\", \"77\u00f73=25, 2\", \"77\u00f73=25, 2\", \"53\u00f74=13, 1\", \"80\u00f77=11, 3\")\n    ,@(\"48\u00f73=16, 0\", \"55\u00f75=11, 0\", \"14\u00f79=1, 5\", \"17\u00f74=4, 1\", \"60\u00f73=20, 0\")\n    ,@(\"30\u00f75=6, 0\", \"68\u00f76=11, 2\", \"61\u00f72=30, 1\", \"79\u00f74=19, 3\", \"69\u00f74=17, 1\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$dataRowIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    if ($dataRowIndex -ge $newValues.Count) {\n        break\n    }\n\n    # Spacer rows have empty paragraphs in every cell; only rows that\n    # actually carry an answer are touched, so this walk is robust even\n    # if the blank-row spacing ever changes.\n    $hasText = $false\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cellText = $t.Cell($r, $c).Range.Text\n        $cellText = $cellText.TrimEnd([char]7).TrimEnd([char]13)\n        if ($cellText.Trim().Length -gt 0) {\n            $hasText = $true\n        }\n    }\n    if (-not $hasText) {\n        continue\n    }\n\n    $rowValues = $newValues[$dataRowIndex]\n    for ($c = 1; $c -le $colCount -and $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n    $dataRowIndex++\n}\n"}
